$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 0.1636683333333333
$ws.Range("H2").Value2 = 0.491005
$ws.Range("I2").Value2 = 0.008639493057305454
$ws.Range("J2").Value2 = 0.008639493057305455
$ws.Range("M2").Value2 = 5.333065333333334
$ws.Range("N2").Value2 = 15.999196
$ws.Range("O2").Value2 = 0.1422335064894009
$ws.Range("P2").Value2 = 0.1422335064894009
$ws.Range("Q2").Value2 = 0.8728539146644445
$ws.Range("R2").Value2 = 7.855685231980001
$ws.Range("S2").Value2 = 0.001228825391831389
$ws.Range("T2").Value2 = 0.00122882539183139

# Row 3
$ws.Range("G3").Value2 = 0.1636683333333333
$ws.Range("H3").Value2 = 0.491005
$ws.Range("I3").Value2 = 0.008639493057305454
$ws.Range("J3").Value2 = 0.008639493057305455
$ws.Range("M3").Value2 = 20.88867166666667
$ws.Range("N3").Value2 = 62.666015
$ws.Range("O3").Value2 = 0.5571034351455781
$ws.Range("P3").Value2 = 0.5571034351455783
$ws.Range("Q3").Value2 = 3.418814077230556
$ws.Range("R3").Value2 = 30.769326695075
$ws.Range("S3").Value2 = 0.004813091260141241
$ws.Range("T3").Value2 = 0.004813091260141243

# Row 4
$ws.Range("G4").Value2 = 0.1636683333333333
$ws.Range("H4").Value2 = 0.491005
$ws.Range("I4").Value2 = 0.008639493057305454
$ws.Range("J4").Value2 = 0.008639493057305455
$ws.Range("M4").Value2 = 11.27340366666667
$ws.Range("N4").Value2 = 33.820211
$ws.Range("O4").Value2 = 0.3006630583650208
$ws.Range("P4").Value2 = 0.3006630583650208
$ws.Range("Q4").Value2 = 1.845099189117222
$ws.Range("R4").Value2 = 16.605892702055
$ws.Range("S4").Value2 = 0.002597576405332822
$ws.Range("T4").Value2 = 0.002597576405332822

# Row 5
$ws.Range("I5").Value2 = 0.808839719627903
$ws.Range("J5").Value2 = 0.8088397196279031
$ws.Range("M5").Value2 = 5.333065333333334
$ws.Range("N5").Value2 = 15.999196
$ws.Range("O5").Value2 = 0.1422335064894009
$ws.Range("P5").Value2 = 0.1422335064894009
$ws.Range("Q5").Value2 = 81.71763214929868
$ws.Range("R5").Value2 = 735.4586893436881
$ws.Range("S5").Value2 = 0.1150441095105806
$ws.Range("T5").Value2 = 0.1150441095105806

# Row 6
$ws.Range("I6").Value2 = 0.808839719627903
$ws.Range("J6").Value2 = 0.8088397196279031
$ws.Range("M6").Value2 = 20.88867166666667
$ws.Range("N6").Value2 = 62.666015
$ws.Range("O6").Value2 = 0.5571034351455781
$ws.Range("P6").Value2 = 0.5571034351455783
$ws.Range("Q6").Value2 = 320.0734813194634
$ws.Range("S6").Value2 = 0.4506073862868911
$ws.Range("T6").Value2 = 0.4506073862868912

# Row 7
$ws.Range("I7").Value2 = 0.808839719627903
$ws.Range("J7").Value2 = 0.8088397196279031
$ws.Range("M7").Value2 = 11.27340366666667
$ws.Range("N7").Value2 = 33.820211
$ws.Range("O7").Value2 = 0.3006630583650208
$ws.Range("P7").Value2 = 0.3006630583650208
$ws.Range("Q7").Value2 = 172.7404028120953
$ws.Range("R7").Value2 = 1554.663625308858
$ws.Range("S7").Value2 = 0.2431882238304313
$ws.Range("T7").Value2 = 0.2431882238304313

# Row 8
$ws.Range("G8").Value2 = 3.457711333333334
$ws.Range("H8").Value2 = 10.373134
$ws.Range("I8").Value2 = 0.1825207873147914
$ws.Range("J8").Value2 = 0.1825207873147914
$ws.Range("M8").Value2 = 5.333065333333334
$ws.Range("N8").Value2 = 15.999196
$ws.Range("O8").Value2 = 0.1422335064894009
$ws.Range("P8").Value2 = 0.1422335064894009
$ws.Range("Q8").Value2 = 18.44020044447378
$ws.Range("R8").Value2 = 165.961804000264
$ws.Range("S8").Value2 = 0.02596057158698895
$ws.Range("T8").Value2 = 0.02596057158698895

# Row 9
$ws.Range("G9").Value2 = 3.457711333333334
$ws.Range("H9").Value2 = 10.373134
$ws.Range("I9").Value2 = 0.1825207873147914
$ws.Range("J9").Value2 = 0.1825207873147914
$ws.Range("M9").Value2 = 20.88867166666667
$ws.Range("N9").Value2 = 62.666015
$ws.Range("O9").Value2 = 0.5571034351455781
$ws.Range("P9").Value2 = 0.5571034351455783
$ws.Range("Q9").Value2 = 72.22699676011223
$ws.Range("R9").Value2 = 650.0429708410101
$ws.Range("S9").Value2 = 0.1016829575985458
$ws.Range("T9").Value2 = 0.1016829575985458

# Row 10
$ws.Range("G10").Value2 = 3.457711333333334
$ws.Range("H10").Value2 = 10.373134
$ws.Range("I10").Value2 = 0.1825207873147914
$ws.Range("J10").Value2 = 0.1825207873147914
$ws.Range("M10").Value2 = 11.27340366666667
$ws.Range("N10").Value2 = 33.820211
$ws.Range("O10").Value2 = 0.3006630583650208
$ws.Range("P10").Value2 = 0.3006630583650208
$ws.Range("Q10").Value2 = 38.98017562347489
$ws.Range("R10").Value2 = 350.821580611274
$ws.Range("S10").Value2 = 0.05487725812925669
$ws.Range("T10").Value2 = 0.05487725812925669
